$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original column A was an empty spacer column (no header, no data) -
# remove it so the "First name / Last name / Occupation" table shifts left
# into A:C (matches dimension going from A1:G21 to A1:F21).
$ws.Range("A1").EntireColumn.Delete()

# Add a new data row under the existing 3 people, this time including a
# genuine numeric value (not a shared-string) to exercise number handling.
$ws.Range("A5").Value = "John"
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = "Tester"

$ws.Range("A6").Select() | Out-Null
